# ITimeU User stories - add priority/estimation values and comments
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Story points ("actual"/estimated value) entered in column C for each user story row.
$ws.Range("C4").Value = 5
$ws.Range("C5").Value = 5
$ws.Range("C6").Value = 3
$ws.Range("C7").Value = 5
$ws.Range("C8").Value = 2
$ws.Range("C9").Value = 4
$ws.Range("C10").Value = 5
$ws.Range("C11").Value = 5
$ws.Range("C12").Value = 2
$ws.Range("C13").Value = 5
$ws.Range("C14").Value = 5
$ws.Range("C15").Value = 3
$ws.Range("C16").Value = 5
$ws.Range("C17").Value = 4
$ws.Range("C18").Value = 4
$ws.Range("C19").Value = 5
$ws.Range("C20").Value = 4
$ws.Range("C21").Value = 1
$ws.Range("C22").Value = 2
$ws.Range("C23").Value = 3
$ws.Range("C24").Value = 5
$ws.Range("C25").Value = 5
$ws.Range("C26").Value = 3
$ws.Range("C27").Value = 3
$ws.Range("C28").Value = 3
$ws.Range("C29").Value = 4
$ws.Range("C30").Value = 1
$ws.Range("C31").Value = 5
$ws.Range("C32").Value = 1
$ws.Range("C34").Value = 2
$ws.Range("C35").Value = 1
$ws.Range("C37").Value = 3
$ws.Range("C38").Value = 2
$ws.Range("C39").Value = 5

# Comments column (E) for a few rows.
$ws.Range("E22").Value = "Could be checked manually I guess"
$ws.Range("E23").Value = "This should normally be part of the import, yes"
$ws.Range("E30").Value = "A warning is enough"

# Extra scoring values added below the main table (not-estimated stories section).
$ws.Range("C48").Value = 5
$ws.Range("C49").Value = 4
$ws.Range("C50").Value = 2

# Move the active selection to A15 (matches the saved view state in the file).
$ws.Range("A15").Select()
